# Temperature Sensor Rev 1 BOM update:
# Remove the "Neodymium Magnets" (McMaster 5862K141) line item from the BOM
# table (worksheet row 13), which shifts every row below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire worksheet row for the Neodymium Magnets BOM line.
# This shifts rows 14:62 up to 13:61 (cell values, styles, and the table's
# structured range all move together).
$ws.Rows(13).Delete()

# The "Item" numbering in column A is independent data entry (not a
# formula), and in the source edit it was left untouched aside from the
# row removal - i.e. it keeps counting 1..11 for the remaining rows
# instead of inheriting the row-13 "8" that the raw row-delete shift
# would otherwise leave behind. Restore the sequential Item numbers for
# the rows that shifted.
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11

# Leave the cursor where the author's last save left it.
$ws.Range("D5").Select() | Out-Null
